$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Refresh simulation results: new random trial data, sorted descending by p (column B).
$ws.Range("A2").Value = 99
$ws.Range("B2").Value = [double]"0.97492344369863104"
$ws.Range("C2").Value = [double]"0.25789701193273901"
$ws.Range("A3").Value = 97
$ws.Range("B3").Value = [double]"0.97422420020359202"
$ws.Range("C3").Value = [double]"0.85539968502074903"
$ws.Range("A4").Value = 36
$ws.Range("B4").Value = [double]"0.94590240255402303"
$ws.Range("C4").Value = [double]"0.28727535177932201"
$ws.Range("A5").Value = 95
$ws.Range("B5").Value = [double]"0.938257587706679"
$ws.Range("C5").Value = [double]"0.35036249099341898"
$ws.Range("A6").Value = 56
$ws.Range("B6").Value = [double]"0.931918684484078"
$ws.Range("C6").Value = [double]"0.39804241554148301"
$ws.Range("A7").Value = 20
$ws.Range("B7").Value = [double]"0.931270998672733"
$ws.Range("C7").Value = [double]"0.82871565479643094"
$ws.Range("A8").Value = 84
$ws.Range("B8").Value = [double]"0.92176562148661501"
$ws.Range("C8").Value = [double]"0.43301052465230599"
$ws.Range("A9").Value = 77
$ws.Range("B9").Value = [double]"0.89881671768124105"
$ws.Range("C9").Value = [double]"0.91556051357701895"
$ws.Range("A10").Value = 32
$ws.Range("B10").Value = [double]"0.82833793417917201"
$ws.Range("C10").Value = [double]"0.93632287594027197"
$ws.Range("A11").Value = 89
$ws.Range("B11").Value = [double]"0.81349507460942305"
$ws.Range("C11").Value = [double]"0.179848870481018"
$ws.Range("A12").Value = 64
$ws.Range("B12").Value = [double]"0.81251708077597595"
$ws.Range("C12").Value = [double]"0.18657992941648499"
$ws.Range("A13").Value = 80
$ws.Range("B13").Value = [double]"0.79785439724937801"
$ws.Range("C13").Value = [double]"0.74265455063774899"
$ws.Range("A14").Value = 17
$ws.Range("B14").Value = [double]"0.79261030206501404"
$ws.Range("C14").Value = [double]"0.94957770108051598"
$ws.Range("A15").Value = 91
$ws.Range("B15").Value = [double]"0.79236264357230601"
$ws.Range("C15").Value = [double]"0.72416091593012299"
$ws.Range("A16").Value = 42
$ws.Range("B16").Value = [double]"0.79153473221170401"
$ws.Range("C16").Value = [double]"0.61542450401362503"
$ws.Range("A17").Value = 98
$ws.Range("B17").Value = [double]"0.76375210997500997"
$ws.Range("C17").Value = [double]"0.54898203033125004"
$ws.Range("A18").Value = 29
$ws.Range("B18").Value = [double]"0.757281704178419"
$ws.Range("C18").Value = [double]"0.88495877963068204"
$ws.Range("A19").Value = 53
$ws.Range("B19").Value = [double]"0.75378441523307105"
$ws.Range("C19").Value = [double]"0.53755819541862404"
$ws.Range("A20").Value = 14
$ws.Range("B20").Value = [double]"0.73598923955831697"
$ws.Range("C20").Value = [double]"0.16273186564763401"
$ws.Range("A21").Value = 58
$ws.Range("B21").Value = [double]"0.73518703179911105"
$ws.Range("C21").Value = [double]"4.6395229592129598E-2"
$ws.Range("A22").Value = 9
$ws.Range("B22").Value = [double]"0.71781631747118302"
$ws.Range("C22").Value = [double]"0.418010229504163"
$ws.Range("A23").Value = 15
$ws.Range("B23").Value = [double]"0.71491115012532003"
$ws.Range("C23").Value = [double]"0.331903672719427"
$ws.Range("A24").Value = 49
$ws.Range("B24").Value = [double]"0.70270876588616205"
$ws.Range("C24").Value = [double]"1.69909910621821E-2"
$ws.Range("A25").Value = 12
$ws.Range("B25").Value = [double]"0.70097148276211296"
$ws.Range("C25").Value = [double]"0.80697005386595699"
$ws.Range("A26").Value = 8
$ws.Range("B26").Value = [double]"0.68679320536493504"
$ws.Range("C26").Value = [double]"0.94880681758879204"
$ws.Range("A27").Value = 74
$ws.Range("B27").Value = [double]"0.67159190043856198"
$ws.Range("C27").Value = [double]"0.42173231578016601"
$ws.Range("A28").Value = 88
$ws.Range("B28").Value = [double]"0.66076887478680701"
$ws.Range("C28").Value = [double]"0.63660184948133502"
$ws.Range("A29").Value = 39
$ws.Range("B29").Value = [double]"0.63494420296593701"
$ws.Range("C29").Value = [double]"0.811772726695146"
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = [double]"0.61923706568047998"
$ws.Range("C30").Value = [double]"0.26691372497690102"
$ws.Range("A31").Value = 19
$ws.Range("B31").Value = [double]"0.61923706568047998"
$ws.Range("C31").Value = [double]"0.26691372497690102"
$ws.Range("A32").Value = 27
$ws.Range("B32").Value = [double]"0.61923706568047998"
$ws.Range("C32").Value = [double]"0.26691372497690102"
$ws.Range("A33").Value = 94
$ws.Range("B33").Value = [double]"0.61923706568047998"
$ws.Range("C33").Value = [double]"0.26691372497690102"
$ws.Range("A34").Value = 57
$ws.Range("B34").Value = [double]"0.60614127712016097"
$ws.Range("C34").Value = [double]"0.69652675863931002"
$ws.Range("A35").Value = 6
$ws.Range("B35").Value = [double]"0.60024200639661496"
$ws.Range("C35").Value = [double]"7.3245617998138501E-2"
$ws.Range("A36").Value = 85
$ws.Range("B36").Value = [double]"0.599362936916089"
$ws.Range("C36").Value = [double]"0.92834609269282897"
$ws.Range("A37").Value = 83
$ws.Range("B37").Value = [double]"0.596964868281185"
$ws.Range("C37").Value = [double]"0.17416469003545099"
$ws.Range("A38").Value = 92
$ws.Range("B38").Value = [double]"0.59337304247759803"
$ws.Range("C38").Value = [double]"0.75245492589179297"
$ws.Range("A39").Value = 24
$ws.Range("B39").Value = [double]"0.58780330025065697"
$ws.Range("C39").Value = [double]"0.93396899609456596"
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = [double]"0.58696803697710698"
$ws.Range("C40").Value = [double]"0.68720126458175201"
$ws.Range("A41").Value = 37
$ws.Range("B41").Value = [double]"0.576867934730242"
$ws.Range("C41").Value = [double]"0.78676015328792104"
$ws.Range("A42").Value = 78
$ws.Range("B42").Value = [double]"0.55965033314593504"
$ws.Range("C42").Value = [double]"0.78698655955783703"
$ws.Range("A43").Value = 33
$ws.Range("B43").Value = [double]"0.55949911671065999"
$ws.Range("C43").Value = [double]"0.99243671472129302"
$ws.Range("A44").Value = 69
$ws.Range("B44").Value = [double]"0.55089526203750605"
$ws.Range("C44").Value = [double]"0.69583286195467897"
$ws.Range("A45").Value = 68
$ws.Range("B45").Value = [double]"0.54581173532131799"
$ws.Range("C45").Value = [double]"0.57033125633842396"
$ws.Range("A46").Value = 60
$ws.Range("B46").Value = [double]"0.50267877809991102"
$ws.Range("C46").Value = [double]"0.52884593508329403"
$ws.Range("A47").Value = 34
$ws.Range("B47").Value = [double]"0.48795119906801598"
$ws.Range("C47").Value = [double]"0.81031639436515901"
$ws.Range("A48").Value = 71
$ws.Range("B48").Value = [double]"0.48191526048989802"
$ws.Range("C48").Value = [double]"6.3064960665017794E-2"
$ws.Range("A49").Value = 51
$ws.Range("B49").Value = [double]"0.47851964762430699"
$ws.Range("C49").Value = [double]"0.53302450459772599"
$ws.Range("A50").Value = 59
$ws.Range("B50").Value = [double]"0.47851964762430699"
$ws.Range("C50").Value = [double]"0.53302450459772599"
$ws.Range("A51").Value = 90
$ws.Range("B51").Value = [double]"0.46446538069166099"
$ws.Range("C51").Value = [double]"0.42829806740164"
$ws.Range("A52").Value = 48
$ws.Range("B52").Value = [double]"0.430012167353374"
$ws.Range("C52").Value = [double]"0.884860324990001"
$ws.Range("A53").Value = 7
$ws.Range("B53").Value = [double]"0.42150366669794498"
$ws.Range("C53").Value = [double]"0.222755918991172"
$ws.Range("A54").Value = 21
$ws.Range("B54").Value = [double]"0.42150366669794498"
$ws.Range("C54").Value = [double]"0.222755918991172"
$ws.Range("A55").Value = 26
$ws.Range("B55").Value = [double]"0.42150366669794498"
$ws.Range("C55").Value = [double]"0.222755918991172"
$ws.Range("A56").Value = 61
$ws.Range("B56").Value = [double]"0.42150366669794498"
$ws.Range("C56").Value = [double]"0.222755918991172"
$ws.Range("A57").Value = 13
$ws.Range("B57").Value = [double]"0.38879122559795398"
$ws.Range("C57").Value = [double]"0.45355088553567302"
$ws.Range("A58").Value = 72
$ws.Range("B58").Value = [double]"0.37644437855021801"
$ws.Range("C58").Value = [double]"0.80977490532521801"
$ws.Range("A59").Value = 43
$ws.Range("B59").Value = [double]"0.37441526147772197"
$ws.Range("C59").Value = [double]"1.9728734871265999E-3"
$ws.Range("A60").Value = 62
$ws.Range("B60").Value = [double]"0.36909321613217"
$ws.Range("C60").Value = [double]"0.94582865395905302"
$ws.Range("A61").Value = 28
$ws.Range("B61").Value = [double]"0.36184962079489902"
$ws.Range("C61").Value = [double]"1.5951822834761899E-2"
$ws.Range("A62").Value = 31
$ws.Range("B62").Value = [double]"0.36184962079489902"
$ws.Range("C62").Value = [double]"1.5951822834761899E-2"
$ws.Range("A63").Value = 44
$ws.Range("B63").Value = [double]"0.36184962079489902"
$ws.Range("C63").Value = [double]"1.5951822834761899E-2"
$ws.Range("A64").Value = 4
$ws.Range("B64").Value = [double]"0.35317473883684902"
$ws.Range("C64").Value = [double]"0.605677549132281"
$ws.Range("A65").Value = 96
$ws.Range("B65").Value = [double]"0.33725842914849202"
$ws.Range("C65").Value = [double]"0.76944558082588999"
$ws.Range("A66").Value = 87
$ws.Range("B66").Value = [double]"0.32966043923215799"
$ws.Range("C66").Value = [double]"0.179670244936292"
$ws.Range("A67").Value = 23
$ws.Range("B67").Value = [double]"0.32937756315966799"
$ws.Range("C67").Value = [double]"0.12875489542062099"
$ws.Range("A68").Value = 70
$ws.Range("B68").Value = [double]"0.32849934894517502"
$ws.Range("C68").Value = [double]"0.417730337267125"
$ws.Range("A69").Value = 65
$ws.Range("B69").Value = [double]"0.32695083644963602"
$ws.Range("C69").Value = [double]"0.46898738146291102"
$ws.Range("A70").Value = 2
$ws.Range("B70").Value = [double]"0.32454977147685798"
$ws.Range("C70").Value = [double]"1.9622903041966001E-2"
$ws.Range("A71").Value = 18
$ws.Range("B71").Value = [double]"0.31131117566531502"
$ws.Range("C71").Value = [double]"2.0756162848575E-2"
$ws.Range("A72").Value = 25
$ws.Range("B72").Value = [double]"0.305065050262884"
$ws.Range("C72").Value = [double]"0.26464738226225898"
$ws.Range("A73").Value = 30
$ws.Range("B73").Value = [double]"0.30104169043092399"
$ws.Range("C73").Value = [double]"0.227658343138399"
$ws.Range("A74").Value = 55
$ws.Range("B74").Value = [double]"0.285576193261875"
$ws.Range("C74").Value = [double]"0.249542249261717"
$ws.Range("A75").Value = 47
$ws.Range("B75").Value = [double]"0.27311053697450199"
$ws.Range("C75").Value = [double]"0.96098124753941505"
$ws.Range("A76").Value = 11
$ws.Range("B76").Value = [double]"0.27216523108002499"
$ws.Range("C76").Value = [double]"0.56746021189958196"
$ws.Range("A77").Value = 67
$ws.Range("B77").Value = [double]"0.26922808056090702"
$ws.Range("C77").Value = [double]"0.72908519542272399"
$ws.Range("A78").Value = 52
$ws.Range("B78").Value = [double]"0.25985428860168203"
$ws.Range("C78").Value = [double]"0.96044720984480303"
$ws.Range("A79").Value = 50
$ws.Range("B79").Value = [double]"0.234324443370894"
$ws.Range("C79").Value = [double]"0.86389825205075199"
$ws.Range("A80").Value = 10
$ws.Range("B80").Value = [double]"0.18686320236439199"
$ws.Range("C80").Value = [double]"0.92313482652757395"
$ws.Range("A81").Value = 5
$ws.Range("B81").Value = [double]"0.17666646878884601"
$ws.Range("C81").Value = [double]"0.18844246436734199"
$ws.Range("A82").Value = 41
$ws.Range("B82").Value = [double]"0.172211683408036"
$ws.Range("C82").Value = [double]"0.81756740846442499"
$ws.Range("A83").Value = 22
$ws.Range("B83").Value = [double]"0.16379535000342599"
$ws.Range("C83").Value = [double]"0.43990456107967202"
$ws.Range("A84").Value = 54
$ws.Range("B84").Value = [double]"0.15001692792260299"
$ws.Range("C84").Value = [double]"0.55149636499285404"
$ws.Range("A85").Value = 81
$ws.Range("B85").Value = [double]"0.134132935900411"
$ws.Range("C85").Value = [double]"0.977784004744785"
$ws.Range("A86").Value = 82
$ws.Range("B86").Value = [double]"0.122572070493049"
$ws.Range("C86").Value = [double]"0.39976280768397998"
$ws.Range("A87").Value = 46
$ws.Range("B87").Value = [double]"0.117020382718869"
$ws.Range("C87").Value = [double]"0.40940448452717199"
$ws.Range("A88").Value = 76
$ws.Range("B88").Value = [double]"0.117020382718869"
$ws.Range("C88").Value = [double]"0.40940448452717199"
$ws.Range("A89").Value = 86
$ws.Range("B89").Value = [double]"0.114918844520482"
$ws.Range("C89").Value = [double]"0.76592030864147598"
$ws.Range("A90").Value = 45
$ws.Range("B90").Value = [double]"9.8624063476444201E-2"
$ws.Range("C90").Value = [double]"0.54585946710943001"
$ws.Range("A91").Value = 16
$ws.Range("B91").Value = [double]"9.1716595007661494E-2"
$ws.Range("C91").Value = [double]"0.84914421244191596"
$ws.Range("A92").Value = 93
$ws.Range("B92").Value = [double]"8.7514057967843506E-2"
$ws.Range("C92").Value = [double]"0.67504053390338403"
$ws.Range("A93").Value = 40
$ws.Range("B93").Value = [double]"8.3168022991189403E-2"
$ws.Range("C93").Value = [double]"0.88767906402437802"
$ws.Range("A94").Value = 79
$ws.Range("B94").Value = [double]"7.0998858178492705E-2"
$ws.Range("C94").Value = [double]"0.55739630908420301"
$ws.Range("A95").Value = 75
$ws.Range("B95").Value = [double]"5.9800329900558599E-2"
$ws.Range("C95").Value = [double]"0.89613977539746104"
$ws.Range("A96").Value = 63
$ws.Range("B96").Value = [double]"5.5365282412294201E-2"
$ws.Range("C96").Value = [double]"0.10127765123867501"
$ws.Range("A97").Value = 0
$ws.Range("B97").Value = [double]"4.6296895400565002E-2"
$ws.Range("C97").Value = [double]"1.1404768794402199E-2"
$ws.Range("A98").Value = 35
$ws.Range("B98").Value = [double]"4.6296895400565002E-2"
$ws.Range("C98").Value = [double]"1.1404768794402199E-2"
$ws.Range("A99").Value = 66
$ws.Range("B99").Value = [double]"4.4269185249727297E-2"
$ws.Range("C99").Value = [double]"0.764303486389598"
$ws.Range("A100").Value = 3
$ws.Range("B100").Value = [double]"1.8331808375969299E-2"
$ws.Range("C100").Value = [double]"0.45297667443607997"
$ws.Range("A101").Value = 73
$ws.Range("B101").Value = [double]"2.7403684536631001E-4"
$ws.Range("C101").Value = [double]"0.820503462613421"

# Player rows that no longer carry the scientific-notation "abstain" highlight.
$ws.Range("B72").ClearFormats()
$ws.Range("B100").ClearFormats()

# Player rows newly flagged (scientific-notation number format) as abstaining.
$ws.Range("B17").NumberFormat = "0.00E+00"
$ws.Range("B68").NumberFormat = "0.00E+00"
$ws.Range("B101").NumberFormat = "0.00E+00"

Write-Output "done"
